$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 2499.8333
$ws.Range("I48").Value = 999.75
$ws.Range("K48").Value = 2999.25
$ws.Range("M48").Value = -2707.25
$ws.Range("H56").Value = 2499.8333
$ws.Range("I56").Value = 999.75
$ws.Range("K56").Value = 2999.25
$ws.Range("M56").Value = -2465.25
$ws.Range("H135").Value = 641.8
$ws.Range("I135").Value = 637.5833
$ws.Range("K135").Value = 5738.2497
$ws.Range("M135").Value = -3203.2497
$ws.Range("H137").Value = 4449.909
$ws.Range("I137").Value = 1463.8125
$ws.Range("K137").Value = 4391.4375
$ws.Range("M137").Value = -1841.4375
$ws.Range("H138").Value = 302189.9
$ws.Range("I138").Value = 6421.8335
$ws.Range("J138").Value = 351484.6
$ws.Range("K138").Value = 19265.5005
$ws.Range("L138").Value = 1054453.8
$ws.Range("M138").Value = -14125.5005
$ws.Range("N138").Value = -1064733.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6295.3877
$ws.Range("I32").Value = 5643.8667
$ws.Range("K32").Value = 5643.8667
$ws.Range("M32").Value = -5356.8667
$ws.Range("H45").Value = 20784.912
$ws.Range("I45").Value = 34321.617
$ws.Range("J45").Value = 3187.2
$ws.Range("K45").Value = 34321.617
$ws.Range("L45").Value = 3187.2
$ws.Range("M45").Value = -33944.617
$ws.Range("N45").Value = -3941.2
$ws.Range("H46").Value = 1393.3334
$ws.Range("J46").Value = 1391
$ws.Range("L46").Value = 1391
$ws.Range("N46").Value = -2029
$ws.Range("H61").Value = 3598.2954
$ws.Range("I61").Value = 1567.9706
$ws.Range("K61").Value = 1567.9706
$ws.Range("M61").Value = -1355.9706
$ws.Range("H130").Value = 88996.5
$ws.Range("J130").Value = 88996.5
$ws.Range("L130").Value = 88996.5
$ws.Range("N130").Value = -99036.5
$ws.Range("H136").Value = 3598.2954
$ws.Range("I136").Value = 1567.9706
$ws.Range("K136").Value = 4703.9118
$ws.Range("M136").Value = -2153.9118

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 6502604
$ws.Range("I105").Value = 359602.4
$ws.Range("K105").Value = 359602.4
$ws.Range("M105").Value = -357855.4
$ws.Range("H132").Value = 99000
$ws.Range("J132").Value = 99000
$ws.Range("L132").Value = 99000
$ws.Range("N132").Value = -109120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1786.9166
$ws.Range("J16").Value = 839.3333
$ws.Range("L16").Value = 839.3333
$ws.Range("N16").Value = -1413.3333
$ws.Range("H17").Value = 35000
$ws.Range("J17").Value = 35000
$ws.Range("L17").Value = 35000
$ws.Range("N17").Value = -35348
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H31").Value = 3443.2297
$ws.Range("I31").Value = 2743.653
$ws.Range("J31").Value = 4814.4
$ws.Range("K31").Value = 2743.653
$ws.Range("L31").Value = 4814.4
$ws.Range("M31").Value = -2448.653
$ws.Range("N31").Value = -5404.4
$ws.Range("H34").Value = 3443.2297
$ws.Range("I34").Value = 2743.653
$ws.Range("J34").Value = 4814.4
$ws.Range("K34").Value = 2743.653
$ws.Range("L34").Value = 4814.4
$ws.Range("M34").Value = -2541.653
$ws.Range("N34").Value = -5218.4
$ws.Range("H58").Value = 2377.7576
$ws.Range("I58").Value = 1347.8636
$ws.Range("J58").Value = 4437.5454
$ws.Range("K58").Value = 1347.8636
$ws.Range("L58").Value = 4437.5454
$ws.Range("M58").Value = -1144.8636
$ws.Range("N58").Value = -4843.5454
$ws.Range("H105").Value = 2509.5
$ws.Range("I105").Value = 1992.125
$ws.Range("J105").Value = 3026.875
$ws.Range("K105").Value = 1992.125
$ws.Range("L105").Value = 3026.875
$ws.Range("M105").Value = -245.125
$ws.Range("N105").Value = -6520.875
$ws.Range("H107").Value = 770
$ws.Range("I107").Value = 874.0909
$ws.Range("J107").Value = 606.4286
$ws.Range("K107").Value = 874.0909
$ws.Range("L107").Value = 606.4286
$ws.Range("M107").Value = 1045.9091
$ws.Range("N107").Value = -4446.4286
$ws.Range("H113").Value = 1786.9166
$ws.Range("J113").Value = 839.3333
$ws.Range("L113").Value = 839.3333
$ws.Range("N113").Value = -5179.3333
$ws.Range("H122").Value = 2304.6667
$ws.Range("I122").Value = 2128.4443
$ws.Range("K122").Value = 6385.3329
$ws.Range("M122").Value = -3935.3329
$ws.Range("H132").Value = 3221.4412
$ws.Range("I132").Value = 2673.4138
$ws.Range("J132").Value = 6400
$ws.Range("K132").Value = 8020.241399999999
$ws.Range("L132").Value = 19200
$ws.Range("M132").Value = -5490.241399999999
$ws.Range("N132").Value = -24260
$ws.Range("H136").Value = 2377.7576
$ws.Range("I136").Value = 1347.8636
$ws.Range("J136").Value = 4437.5454
$ws.Range("K136").Value = 4043.5908
$ws.Range("L136").Value = 13312.6362
$ws.Range("M136").Value = -1493.5908
$ws.Range("N136").Value = -18412.6362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1026.7826
$ws.Range("I2").Value = 66.85714
$ws.Range("K2").Value = 401.14284
$ws.Range("M2").Value = -288.14284
$ws.Range("H104").Value = 3749.5
$ws.Range("J104").Value = 6500
$ws.Range("L104").Value = 19500
$ws.Range("N104").Value = -24742
$ws.Range("H107").Value = 828.8823
$ws.Range("J107").Value = 872.93335
$ws.Range("L107").Value = 2618.80005
$ws.Range("N107").Value = -6458.80005
$ws.Range("H113").Value = 1570.6364
$ws.Range("J113").Value = 1918.4
$ws.Range("L113").Value = 5755.200000000001
$ws.Range("N113").Value = -10095.2
$ws.Range("H122").Value = 1405.091
$ws.Range("I122").Value = 1318
$ws.Range("J122").Value = 1418.8422
$ws.Range("K122").Value = 11862
$ws.Range("L122").Value = 12769.5798
$ws.Range("M122").Value = -9412
$ws.Range("N122").Value = -17669.5798
$ws.Range("H131").Value = 2509.0833
$ws.Range("J131").Value = 1922.375
$ws.Range("L131").Value = 5767.125
$ws.Range("N131").Value = -15847.125
$ws.Range("H137").Value = 5267.4375
$ws.Range("J137").Value = 4536.8184
$ws.Range("L137").Value = 13610.4552
$ws.Range("N137").Value = -23810.4552
$ws.Range("H138").Value = 5982.2856
$ws.Range("I138").Value = 4167.7144
$ws.Range("J138").Value = 7796.857
$ws.Range("K138").Value = 12503.1432
$ws.Range("L138").Value = 23390.571
$ws.Range("M138").Value = -7363.143199999999
$ws.Range("N138").Value = -33670.571
$ws.Range("H139").Value = 3122.9
$ws.Range("I139").Value = 2245.8
$ws.Range("K139").Value = 6737.400000000001
$ws.Range("M139").Value = -1597.400000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1888.92
$ws.Range("J102").Value = 2901.1667
$ws.Range("L102").Value = 2901.1667
$ws.Range("N102").Value = -6145.1667
$ws.Range("H107").Value = 288.25
$ws.Range("I107").Value = 132.36363
$ws.Range("K107").Value = 132.36363
$ws.Range("M107").Value = 1787.63637
$ws.Range("H133").Value = 178799.4
$ws.Range("J133").Value = 178799.4
$ws.Range("L133").Value = 178799.4
$ws.Range("N133").Value = -188919.4
$ws.Range("H135").Value = 69065.47
$ws.Range("J135").Value = 69065.47
$ws.Range("L135").Value = 69065.47
$ws.Range("N135").Value = -79205.47
$ws.Range("H140").Value = 74999.234
$ws.Range("J140").Value = 74999.234
$ws.Range("L140").Value = 74999.234
$ws.Range("N140").Value = -85359.234

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("H46").Value = 1547.3948
$ws.Range("J46").Value = 1428.3667
$ws.Range("L46").Value = 1428.3667
$ws.Range("N46").Value = -1804.3667
$ws.Range("H132").Value = 5972.737
$ws.Range("I132").Value = 3375.5
$ws.Range("J132").Value = 6665.3335
$ws.Range("K132").Value = 10126.5
$ws.Range("L132").Value = 19996.0005
$ws.Range("M132").Value = -7596.5
$ws.Range("N132").Value = -25056.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 15999
$ws.Range("J20").Value = 15999
$ws.Range("L20").Value = 15999
$ws.Range("N20").Value = -16479
$ws.Range("H113").Value = 1327.0714
$ws.Range("J113").Value = 1066.6666
$ws.Range("L113").Value = 3199.9998
$ws.Range("N113").Value = -7539.9998
$ws.Range("H136").Value = 28574088
$ws.Range("I136").Value = 31251220
$ws.Range("K136").Value = 93753660
$ws.Range("M136").Value = -93751110
